# Adds a new QoQ forecast column (BB) mirroring column BA, with updated
# tail values for the most recent quarters, plus a new trailing row (83)
# that extends the date axis by one more quarter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header date (column BB, row 1) - copy style from BA1 (date format)
$ws.Range("BB1").Value = 45986
$ws.Range("BA1").Copy() | Out-Null
$ws.Range("BB1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Column BB mirrors column BA for rows 2-70
for ($r = 2; $r -le 70; $r++) {
    $ws.Cells.Item($r, 54).Value = $ws.Cells.Item($r, 53).Value2
}

# Rows 71-82 get revised tail values in column BB
$bbTail = @{
    71 = -0.5
    72 = 0.4
    73 = 0.2
    74 = 0.2
    75 = 0.2
    76 = 0.2
    77 = 0.2
    78 = 0.2
    79 = 0.2
    80 = 0.2
    81 = 0.2
    82 = 0.2
}
foreach ($r in $bbTail.Keys) {
    $ws.Cells.Item($r, 54).Value = $bbTail[$r]
}

# New row 83: extends date column A and adds the BB value
$ws.Range("A83").Value = 46934
$ws.Range("A82").Copy() | Out-Null
$ws.Range("A83").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("BB83").Value = 0.2
